$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.415.97"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.074.89"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.93"
$ws.Range("E5").Value = "  +8.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.54"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("E7").Value = "  -6.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.361"
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.074.23"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.708"
$ws.Range("E11").Value = "  -6.14%  "
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.01"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.329.21"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.34"
$ws.Range("E16").Value = "  -6.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.643.74"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.082.52"
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000210"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.71"
$ws.Range("E21").Value = "  -5.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.95"
$ws.Range("E22").Value = "  -8.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.36"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.70"
$ws.Range("E24").Value = "  -4.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.71"
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.53"
$ws.Range("E26").Value = "  -9.09%  "
$ws.Range("E27").Value = "  -5.62%  "
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("E30").Value = "  +14.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.00"
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("E32").Value = "  -4.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.195"
$ws.Range("E33").Value = "  -17.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.49"
$ws.Range("E34").Value = "  -6.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.150"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.07"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "491.38"
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("E40").Value = "  -5.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("E41").Value = "  +54.09%  "
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.06"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("E45").Value = "  -7.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.73"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("E47").Value = "  -7.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.668"
$ws.Range("E48").Value = "  -8.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.37"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.29"
$ws.Range("E51").Value = "  -5.44%  "
